# Update NATMI edge-weight metrics for Wnt4-Fzd6 LR pair (Young D4) per
# Dr Hou's advice: ligand- and receptor-expressing cell counts change
# from 1 to 3 for every sending/target cluster combination, which in
# turn changes the derived total expression values, specificities and
# edge weights for rows 2-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.408675
$ws.Range("H2").Value = 4.226025
$ws.Range("I2").Value = 0.2345535261744868
$ws.Range("J2").Value = 0.2345535261744868
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.545366
$ws.Range("N2").Value = 61.636098
$ws.Range("O2").Value = 0.8965950288338865
$ws.Range("P2").Value = 0.8965950288338865
$ws.Range("Q2").Value = 28.94174345005
$ws.Range("R2").Value = 260.47569105045
$ws.Range("S2").Value = 0.2102995255635037
$ws.Range("T2").Value = 0.2102995255635038
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.408675
$ws.Range("H3").Value = 4.226025
$ws.Range("I3").Value = 0.2345535261744868
$ws.Range("J3").Value = 0.2345535261744868
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.154739666666666
$ws.Range("N3").Value = 6.464219
$ws.Range("O3").Value = 0.09403234157836461
$ws.Range("P3").Value = 0.09403234157836463
$ws.Range("Q3").Value = 3.035327899941666
$ws.Range("R3").Value = 27.317951099475
$ws.Range("S3").Value = 0.02205561729164923
$ws.Range("T3").Value = 0.02205561729164923
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.408675
$ws.Range("H4").Value = 4.226025
$ws.Range("I4").Value = 0.2345535261744868
$ws.Range("J4").Value = 0.2345535261744868
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2147726666666667
$ws.Range("N4").Value = 0.6443180000000001
$ws.Range("O4").Value = 0.009372629587748921
$ws.Range("P4").Value = 0.009372629587748921
$ws.Range("Q4").Value = 0.3025448862166667
$ws.Range("R4").Value = 2.72290397595
$ws.Range("S4").Value = 0.002198383319333836
$ws.Range("T4").Value = 0.002198383319333836
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.436364333333333
$ws.Range("H5").Value = 4.309093
$ws.Range("I5").Value = 0.2391639798069812
$ws.Range("J5").Value = 0.2391639798069812
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 20.545366
$ws.Range("N5").Value = 61.636098
$ws.Range("O5").Value = 0.8965950288338865
$ws.Range("P5").Value = 0.8965950288338865
$ws.Range("Q5").Value = 29.51063093767933
$ws.Range("R5").Value = 265.595678439114
$ws.Range("S5").Value = 0.2144332353710674
$ws.Range("T5").Value = 0.2144332353710674
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.436364333333333
$ws.Range("H6").Value = 4.309093
$ws.Range("I6").Value = 0.2391639798069812
$ws.Range("J6").Value = 0.2391639798069812
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.154739666666666
$ws.Range("N6").Value = 6.464219
$ws.Range("O6").Value = 0.09403234157836461
$ws.Range("P6").Value = 0.09403234157836463
$ws.Range("Q6").Value = 3.094991204818555
$ws.Range("R6").Value = 27.854920843367
$ws.Range("S6").Value = 0.02248914904245116
$ws.Range("T6").Value = 0.02248914904245116
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.436364333333333
$ws.Range("H7").Value = 4.309093
$ws.Range("I7").Value = 0.2391639798069812
$ws.Range("J7").Value = 0.2391639798069812
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.2147726666666667
$ws.Range("N7").Value = 0.6443180000000001
$ws.Range("O7").Value = 0.009372629587748921
$ws.Range("P7").Value = 0.009372629587748921
$ws.Range("Q7").Value = 0.3084917981748889
$ws.Range("R7").Value = 2.776426183574
$ws.Range("S7").Value = 0.002241595393462698
$ws.Range("T7").Value = 0.002241595393462698
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.397773
$ws.Range("H8").Value = 4.193319
$ws.Range("I8").Value = 0.2327382724485711
$ws.Range("J8").Value = 0.2327382724485711
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 20.545366
$ws.Range("N8").Value = 61.636098
$ws.Range("O8").Value = 0.8965950288338865
$ws.Range("P8").Value = 0.8965950288338865
$ws.Range("Q8").Value = 28.717757869918
$ws.Range("R8").Value = 258.459820829262
$ws.Range("S8").Value = 0.2086719780967756
$ws.Range("T8").Value = 0.2086719780967756
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.397773
$ws.Range("H9").Value = 4.193319
$ws.Range("I9").Value = 0.2327382724485711
$ws.Range("J9").Value = 0.2327382724485711
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.154739666666666
$ws.Range("N9").Value = 6.464219
$ws.Range("O9").Value = 0.09403234157836461
$ws.Range("P9").Value = 0.09403234157836463
$ws.Range("Q9").Value = 3.011836928095666
$ws.Range("R9").Value = 27.106532352861
$ws.Range("S9").Value = 0.02188492473324253
$ws.Range("T9").Value = 0.02188492473324253
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.397773
$ws.Range("H10").Value = 4.193319
$ws.Range("I10").Value = 0.2327382724485711
$ws.Range("J10").Value = 0.2327382724485711
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.2147726666666667
$ws.Range("N10").Value = 0.6443180000000001
$ws.Range("O10").Value = 0.009372629587748921
$ws.Range("P10").Value = 0.009372629587748921
$ws.Range("Q10").Value = 0.3002034346046667
$ws.Range("R10").Value = 2.701830911442
$ws.Range("S10").Value = 0.002181369618553047
$ws.Range("T10").Value = 0.002181369618553047
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.34152
$ws.Range("H11").Value = 4.02456
$ws.Range("I11").Value = 0.2233717830114097
$ws.Range("J11").Value = 0.2233717830114097
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 20.545366
$ws.Range("N11").Value = 61.636098
$ws.Range("O11").Value = 0.8965950288338865
$ws.Range("P11").Value = 0.8965950288338865
$ws.Range("Q11").Value = 27.56201939632
$ws.Range("R11").Value = 248.05817456688
$ws.Range("S11").Value = 0.2002740302297915
$ws.Range("T11").Value = 0.2002740302297915
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.34152
$ws.Range("H12").Value = 4.02456
$ws.Range("I12").Value = 0.2233717830114097
$ws.Range("J12").Value = 0.2233717830114097
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.154739666666666
$ws.Range("N12").Value = 6.464219
$ws.Range("O12").Value = 0.09403234157836461
$ws.Range("P12").Value = 0.09403234157836463
$ws.Range("Q12").Value = 2.890626357626667
$ws.Range("R12").Value = 26.01563721864
$ws.Range("S12").Value = 0.02100417179909722
$ws.Range("T12").Value = 0.02100417179909722
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.34152
$ws.Range("H13").Value = 4.02456
$ws.Range("I13").Value = 0.2233717830114097
$ws.Range("J13").Value = 0.2233717830114097
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.2147726666666667
$ws.Range("N13").Value = 0.6443180000000001
$ws.Range("O13").Value = 0.009372629587748921
$ws.Range("P13").Value = 0.009372629587748921
$ws.Range("Q13").Value = 0.2881218277866667
$ws.Range("R13").Value = 2.59309645008
$ws.Range("S13").Value = 0.00209358098252097
$ws.Range("T13").Value = 0.00209358098252097
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4214396666666667
$ws.Range("H14").Value = 1.264319
$ws.Range("I14").Value = 0.07017243855855113
$ws.Range("J14").Value = 0.07017243855855113
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 20.545366
$ws.Range("N14").Value = 61.636098
$ws.Range("O14").Value = 0.8965950288338865
$ws.Range("P14").Value = 0.8965950288338865
$ws.Range("Q14").Value = 8.658632198584668
$ws.Range("R14").Value = 77.92768978726201
$ws.Range("S14").Value = 0.06291625957274828
$ws.Range("T14").Value = 0.06291625957274828
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4214396666666667
$ws.Range("H15").Value = 1.264319
$ws.Range("I15").Value = 0.07017243855855113
$ws.Range("J15").Value = 0.07017243855855113
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.154739666666666
$ws.Range("N15").Value = 6.464219
$ws.Range("O15").Value = 0.09403234157836461
$ws.Range("P15").Value = 0.09403234157836463
$ws.Range("Q15").Value = 0.9080927668734444
$ws.Range("R15").Value = 8.172834901861002
$ws.Range("S15").Value = 0.006598478711924483
$ws.Range("T15").Value = 0.006598478711924484
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4214396666666667
$ws.Range("H16").Value = 1.264319
$ws.Range("I16").Value = 0.07017243855855113
$ws.Range("J16").Value = 0.07017243855855113
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.2147726666666667
$ws.Range("N16").Value = 0.6443180000000001
$ws.Range("O16").Value = 0.009372629587748921
$ws.Range("P16").Value = 0.009372629587748921
$ws.Range("Q16").Value = 0.09051372104911114
$ws.Range("R16").Value = 0.8146234894420002
$ws.Range("S16").Value = 0.0006577002738783695
$ws.Range("T16").Value = 0.0006577002738783695
